# Append the new "2025-08" stats row (row 21) to the sheet, matching the
# existing table layout (columns A-H: month, schools, authorities, users,
# users_per_school, yoy_schools, yoy_authorities, yoy_users).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 21

$ws.Range("A$newRow").Value = 45870
$ws.Range("B$newRow").Value = 6208
$ws.Range("C$newRow").Value = 980
$ws.Range("D$newRow").Value = 5584698
$ws.Range("E$newRow").Value = 899.5969716494845
$ws.Range("F$newRow").Value = 7.759069605971192
$ws.Range("G$newRow").Value = 3.375527426160341
$ws.Range("H$newRow").Value = 27.4762313020946

# Match the date/time number format used by the rest of column A.
$ws.Range("A$newRow").NumberFormat = $ws.Range("A20").NumberFormat
